$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Spanish word translations (column D) to infinitive form
$ws.Range("D2").Value = "bañarse"
$ws.Range("D3").Value = "abotonar"

# Update "correct" answer column (G) from single-letter key codes to left/right labels
$ws.Range("G2").Value = "left"
$ws.Range("G3").Value = "left"
$ws.Range("G4").Value = "right"
$ws.Range("G5").Value = "right"

# The "correct" column no longer uses the centered style - reset to General/default style
$ws.Range("G2:G5").Style = "Normal"

# Column width adjustments: D:E were auto-fit to the new (shorter) contents,
# F was manually widened to fit the longer "POS"/file path text
$ws.Columns.Item(4).ColumnWidth = 7.91796875
$ws.Columns.Item(5).ColumnWidth = 3.91796875
$ws.Columns.Item(6).ColumnWidth = 35.91796875

# Update the active cell selection
$ws.Range("F11").Select()
